$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("module3")

# Row 13 - update the "fim" (end) time of the existing lesson; start (C13) stays the same
$ws.Range("D13").Value = 0.9458333333333333

# Row 14 - new lesson: Capitulo 19 Aula 02
$ws.Range("A14").Value = "Capítulo 19 Aula 02 – Colocando uma imagem de fundo no seu site"
$ws.Range("B14").Value = 15
$ws.Range("C14").Value = 0.9458333333333333
$ws.Range("D14").Value = 0.9569444444444444
$ws.Range("E14").Formula = "=D14-C14"
$ws.Range("F14").Value = 44980
$ws.Range("G14").Formula = "=SUM(E2:E14)+module2!G51"

# Row 15 - new lesson: Capitulo 19 Aula 03
$ws.Range("A15").Value = "Capítulo 19 Aula 03 – Imagens que se repetem no fundo do site"
$ws.Range("B15").Value = 10
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0.0062499999999999995
$ws.Range("E15").Formula = "=D15-C15"
$ws.Range("F15").Value = 44981

# Row 16 - new lesson: Capitulo 19 Aula 04
$ws.Range("A16").Value = "Capítulo 19 Aula 04 – Configurando a posição da imagem no fundo do site"
$ws.Range("B16").Value = 17
$ws.Range("C16").Value = 0.007638888888888889
$ws.Range("D16").Value = 0.035416666666666666
$ws.Range("E16").Formula = "=D16-C16"
$ws.Range("F16").Value = 44981

# Row 17 - new lesson: Capitulo 19 Aula 05
$ws.Range("A17").Value = "Capítulo 19 Aula 05 – Mudando o tamanho da imagem de fundo do site"
$ws.Range("B17").Value = 16
$ws.Range("C17").Value = 0.5583333333333333
$ws.Range("D17").Value = 0.5750000000000001
$ws.Range("E17").Formula = "=D17-C17"
$ws.Range("F17").Value = 44981

# Row 18 - new lesson: Capitulo 19 Aula 06
$ws.Range("A18").Value = "Capítulo 19 Aula 06 – background-attachment e shorthand"
$ws.Range("B18").Value = 11
$ws.Range("C18").Value = 0.5756944444444444
$ws.Range("D18").Value = 0.5875
$ws.Range("E18").Formula = "=D18-C18"
$ws.Range("F18").Value = 44981
# the running-total formula that used to live on G18 now belongs on G14 (see above),
# so fully remove this cell (not just its contents)
$ws.Range("G18").Clear()

# Rows 19-20 get the same "day" marker as the new lessons above
$ws.Range("F19").Value = 44981
$ws.Range("F20").Value = 44981

# Move the current selection to reflect where the user ended up after editing
$ws.Range("A20").Select()

$wb.Save()
